# Diary workbook: add last week's entries (rows 27-30) that were
# previously blank placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: 2020-02-13 (serial 43874) ---
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)   # xlPasteFormats - reuse existing date style
$ws.Range("A27").Value = 43874
$ws.Range("B27").Value = "09:00 am - 12:00 pm"
$ws.Range("C27").Value = "N/A"
$ws.Range("D27").Value = "Complete revision for the test"
$ws.Range("E27").Value = "Referred to the slides and some articles on the internet to get a better understanding of the concepts like all the different types of structural and behavioral diagrams."
$ws.Range("F27").Value = "Since I had paid attention in class, it was easier to revise, as after reiterating the examples given in class to explain the concepts my understanding of it had greatly increased. "
$ws.Range("G27").Value = "I feel like Ive prepared enough. But still nervous as this is my first exam in 3 years!"
$ws.Rows.Item(27).RowHeight = 93.6

# --- Row 28: 2020-02-13 (serial 43874) ---
$ws.Range("A26").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 43874
$ws.Range("B28").Value = "05:00 pm - 07:50 pm"
$ws.Range("C28").Value = "N/A"
$ws.Range("D28").Value = "Successfully complete the test and learn something new!"
$ws.Range("E28").Value = "Test went well! Also learned other points of KEP and the big picture like the stakeholders, the developers etc."
$ws.Range("F28").Value = "I feel like I wasted my time on answers when it could've been completed in less time and that`u{2019}s why I had to rush a little in the end and couldn't really write down all the points I wanted to for the last question. But overall good paper"
$ws.Range("G28").Value = "Next time I will manage my time better while giving an exam."
$ws.Rows.Item(28).RowHeight = 109.2

# --- Row 29: 2020-02-15 (serial 43876) ---
$ws.Range("A26").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 43876
$ws.Range("B29").Value = "11:00 am - 02:00 pm"
$ws.Range("C29").Value = "Anjana, Aman"
$ws.Range("D29").Value = "To resubmit our first homework and start with the next homework"
$ws.Range("E29").Value = "We added the missing UML diagrams and also explained the flow more clearly this time.Explained the features of the application."
$ws.Range("F29").Value = "Gained a better understanding of how explanation of the flow must be done. Like giving more diagrams so the reader can relate to the given explanation of it better. "
$ws.Range("G29").Value = "Hopefully our score for the frist assignment would increase."
$ws.Rows.Item(29).RowHeight = 78

# --- Row 30: 2020-02-19 (serial 43880) ---
$ws.Range("A26").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = 43880
$ws.Range("B30").Value = "11:00 pm - 12:00 am"
$ws.Range("C30").Value = "Anjana, Aman"
$ws.Range("D30").Value = "To complete our assigment 3."
$ws.Range("E30").Value = "We completed our assignment 3 and put as much as relevant information needed. Formatted the document so that it would be pleasing to the reader.Succesfully selected the open issues which we think we can tackle. Also managed to find stakeholders of the application."
$ws.Range("F30").Value = "It had many stakeholders to begin with. Realised why our application is unique. "
$ws.Range("G30").Value = "Hopefully we have done a better job and that we don`u{2019}t have to resubmit this assigment too."
$ws.Rows.Item(30).RowHeight = 140.4

# --- Update the saved view state (user scrolled / selected E31 last) ---
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 29
$ws.Range("E31").Select()

Write-Output "diary rows 27-30 populated"
